# The commit removes the final slide ("Zoom Poll", sldId 283 / rId33)
# along with its notes page. Delete the last slide in the deck.
$p = $ppt.ActivePresentation
$lastIndex = $p.Slides.Count
$p.Slides.Item($lastIndex).Delete()
